$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeAddr, $val) {
    $rng = $ws.Range($rangeAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

Set-TextValue "D2" '60.162.06'
Set-TextValue "E2" '  +3.73%  '
Set-TextValue "D3" '2.422.64'
Set-TextValue "E3" '  +3.34%  '
Set-TextValue "D5" '554.69'
Set-TextValue "E5" '  +2.38%  '
Set-TextValue "D6" '137.65'
Set-TextValue "E6" '  +2.59%  '
Set-TextValue "E7" '  -0.02%  '
Set-TextValue "D8" '0.578'
Set-TextValue "E8" '  +2.53%  '
Set-TextValue "E9" '  +3.65%  '
Set-TextValue "E10" '  +4.42%  '
Set-TextValue "E11" '  +1.08%  '
Set-TextValue "D13" '24.93'
Set-TextValue "E13" '  +4.75%  '
Set-TextValue "D14" '2.850.93'
Set-TextValue "E14" '  +3.28%  '
Set-TextValue "D15" '60.040.44'
Set-TextValue "E15" '  +3.64%  '
Set-TextValue "E16" '  +3.10%  '
Set-TextValue "D17" '2.419.60'
Set-TextValue "E17" '  +3.03%  '
Set-TextValue "D18" '11.35'
Set-TextValue "E18" '  +6.39%  '
Set-TextValue "D19" '4.39'
Set-TextValue "E19" '  +2.27%  '
Set-TextValue "D20" '332.21'
Set-TextValue "D21" '6.74'
Set-TextValue "E21" '  -0.20%  '
Set-TextValue "E22" '  +0.11%  '
Set-TextValue "D23" '65.39'
Set-TextValue "E23" '  +3.97%  '
Set-TextValue "E24" '  +3.52%  '
Set-TextValue "D25" '8.61'
Set-TextValue "E25" '  +3.56%  '
Set-TextValue "E26" '  +0.44%  '
Set-TextValue "E27" '  +1.05%  '
Set-TextValue "D28" '0.0₃0785'
Set-TextValue "E28" '  +6.71%  '
Set-TextValue "E29" '  +1.20%  '
Set-TextValue "D30" '169.83'
Set-TextValue "E30" '  -0.30%  '
Set-TextValue "E31" '  +1.85%  '
Set-TextValue "B32" 'EthereumClassic'
Set-TextValue "C32" 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue "D32" '18.67'
Set-TextValue "E32" '  +1.90%  '
Set-TextValue "B33" 'SuiNetwork'
Set-TextValue "C33" 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
Set-TextValue "D33" '1.03'
Set-TextValue "E33" '  +2.19%  '
Set-TextValue "E34" '  +0.00%  '
Set-TextValue "E35" '  +5.32%  '
Set-TextValue "E36" '  +0.11%  '
Set-TextValue "D37" '4.22'
Set-TextValue "E37" '  +1.57%  '
Set-TextValue "D38" '1.62'
Set-TextValue "E38" '  +0.78%  '
Set-TextValue "D39" '39.55'
Set-TextValue "E39" '  +1.24%  '
Set-TextValue "D40" '0.417'
Set-TextValue "E40" '  +10.71%  '
Set-TextValue "D41" '313.59'
Set-TextValue "E41" '  +8.59%  '
Set-TextValue "E42" '  +1.64%  '
Set-TextValue "D43" '139.09'
Set-TextValue "E43" '  -1.65%  '
Set-TextValue "D44" '0.0961'
Set-TextValue "E44" '  +1.67%  '
Set-TextValue "D46" '19.50'
Set-TextValue "E46" '  +2.76%  '
Set-TextValue "E47" '  +8.89%  '
Set-TextValue "D48" '0.576'
Set-TextValue "E48" '  +1.70%  '
Set-TextValue "E49" '  +1.64%  '
Set-TextValue "D50" '17.72'
Set-TextValue "E50" '  +2.01%  '
Set-TextValue "D51" '11.06'
Set-TextValue "E51" '  -0.11%  '
